$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.698.05"
$ws.Range("E2").Value = "  -3.30%  "

# Row 3
$ws.Range("D3").Value = "3.653.50"
$ws.Range("E3").Value = "  +2.64%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.20%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "405.77"
$ws.Range("E5").Value = "  -2.83%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.60"
$ws.Range("E6").Value = "  +1.49%  "

# Row 7
$ws.Range("D7").Value = "3.649.11"
$ws.Range("E7").Value = "  +2.80%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.622"
$ws.Range("E8").Value = "  -4.37%  "

# Row 9
$ws.Range("E9").Value = "  +0.20%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.721"
$ws.Range("E10").Value = "  -7.79%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.161"
$ws.Range("E11").Value = "  -10.98%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000313"
$ws.Range("E12").Value = "  -11.03%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.87"
$ws.Range("E13").Value = "  -1.89%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.87"
$ws.Range("E14").Value = "  -0.86%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "4.230.39"
$ws.Range("E15").Value = "  +2.78%  "

# Row 17
$ws.Range("D17").Value = "3.648.73"
$ws.Range("E17").Value = "  +3.25%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.82"
$ws.Range("E18").Value = "  -1.93%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.38"
$ws.Range("E19").Value = "  +7.68%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.08"
$ws.Range("E20").Value = "  -4.57%  "

# Row 21
$ws.Range("D21").Value = "64.787.19"
$ws.Range("E21").Value = "  -2.89%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "419.22"
$ws.Range("E22").Value = "  -8.20%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.23"
$ws.Range("E23").Value = "  +16.88%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.68"
$ws.Range("E24").Value = "  -4.61%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.99"
$ws.Range("E25").Value = "  -5.85%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "35.79"
$ws.Range("E26").Value = "  +3.12%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.17"
$ws.Range("E27").Value = "  -6.14%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.38"
$ws.Range("E28").Value = "  -5.75%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.03"
$ws.Range("E29").Value = "  +4.10%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.45"
$ws.Range("E30").Value = "  +0.32%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.70"
$ws.Range("E31").Value = "  -2.93%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.118"
$ws.Range("E32").Value = "  +0.60%  "

# Row 33
$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.92"
$ws.Range("E33").Value = "  -4.80%  "

# Row 34
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.159"
$ws.Range("E34").Value = "  +0.19%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "40.33"
$ws.Range("E35").Value = "  +2.99%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.74"
$ws.Range("E36").Value = "  -1.82%  "

# Row 37
$ws.Range("E37").Value = "  -0.05%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0462"
$ws.Range("E38").Value = "  -6.46%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.91"
$ws.Range("E39").Value = "  +27.06%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.994"
$ws.Range("E40").Value = "  -0.50%  "

# Row 41
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.139"
$ws.Range("E41").Value = "  -5.91%  "

# Row 42
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.86"
$ws.Range("E42").Value = "  +25.38%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.31"
$ws.Range("E43").Value = "  +2.63%  "

# Row 44
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.0₃0628"
$ws.Range("E44").Value = "  -20.24%  "

# Row 45
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.34"
$ws.Range("E45").Value = "  -0.24%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.06"
$ws.Range("E46").Value = "  +17.95%  "

# Row 47
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.34"
$ws.Range("E47").Value = "  -3.82%  "

# Row 48
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.05"
$ws.Range("E48").Value = "  +4.31%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.79"
$ws.Range("E49").Value = "  -5.38%  "

# Row 50
$ws.Range("E50").Value = "  -7.99%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.289"
$ws.Range("E51").Value = "  -5.70%  "
